# Append a new row (row 78) of log data to each of the 4 worksheets,
# matching the next day's reading that follows the existing row 77 data.

$wb = $excel.ActiveWorkbook

# sheetName -> values for columns A..I of the new row
$newRowsBySheet = @{
    "MID_LFT_#1" = @(45864.46186342592, "0x01,0x90", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,", "0x01,0x38", "0x07", 400, [double]"5.68631262647113e+23", 312, 7)
    "MID_LFT_#2" = @(45864.46186342592, "0x01,0x7c", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,", "0x01,0x34", "0x19", 380, [double]"5.68432987514711e+23", 308, 25)
    "MID_PLT_#1" = @(45864.46186342592, "0x00,0x6e", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,", "0x00,0x61", "0x15", 110, [double]"5.68631262647113e+23", 97, 15)
    "MID_PLT_#2" = @(45864.46186342592, "0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,", "0x00,0x77", "0x9", 130, [double]"5.68631262647113e+23", 119, 9)
}

foreach ($sheetName in $newRowsBySheet.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $newRowsBySheet[$sheetName]

    $newRow = $ws.Range("A77").Row + 1

    $ws.Cells.Item($newRow, 1).Value = $values[0]
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

    $ws.Cells.Item($newRow, 2).Value = $values[1]
    $ws.Cells.Item($newRow, 3).Value = $values[2]
    $ws.Cells.Item($newRow, 4).Value = $values[3]
    $ws.Cells.Item($newRow, 5).Value = $values[4]

    $ws.Cells.Item($newRow, 6).Value = $values[5]
    $ws.Cells.Item($newRow, 7).Value = $values[6]
    $ws.Cells.Item($newRow, 8).Value = $values[7]
    $ws.Cells.Item($newRow, 9).Value = $values[8]
}
